$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.390.54"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.562.88"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'286.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").Value = "'0.3652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("D8").Value = "'49.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "'0.3347"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'1.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D11").Value = "'0.07389"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'20.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "'5.919"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "'6.867"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "'88.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").Value = "'0.06748"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'6.293"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'16.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "'11.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "22.370.90"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'2.374"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "'2.532"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").Value = "'149.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "'19.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").Value = "'4.992"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'123.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "1.738.77"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "'6.080"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "'1.989"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'9.533"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("D36").Value = "'0.08247"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "'0.02383"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("E38").Value = "  -5.34%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2216"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("D41").Value = "'5.307"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").Value = "'11.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("D43").Value = "'0.6053"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'13.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").Value = "'3.764"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'0.5727"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").Value = "'2.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("D49").Value = "'124.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("D50").Value = "'1.212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'0.07227"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
